# Append a new data row (surveillance week 47) to the flu-cases table on
# Sheet1, formatted the same way as the existing data rows, and leave the
# selection where it ended up after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 14
$lastRow = $newRow - 1

# Clone the formatting (font, fill, borders, alignment) of the row above
# onto the new row before putting the data in.
$ws.Range("A" + $lastRow + ":G" + $lastRow).Copy()
$ws.Range("A" + $newRow + ":G" + $newRow).PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Rows.Item($newRow).RowHeight = 18

# Week 47 figures.
$ws.Cells.Item($newRow, 1).Value = 47
$ws.Cells.Item($newRow, 2).Value = 5424
$ws.Cells.Item($newRow, 3).Value = 2214
$ws.Cells.Item($newRow, 4).Value = 119
$ws.Cells.Item($newRow, 5).Value = 16
$ws.Cells.Item($newRow, 6).Value = 23.1
$ws.Cells.Item($newRow, 7).Value = 0

# Selection ends up at C18 after the edit.
$ws.Range("C18").Select()
